{"js": "// Ver 7 - Mid Level Contact\n//\n// The commit being replayed:\n//   1) Drops the stray <w:proofErr gramStart/.../gramEnd/> markers that\n//      bracket the \"Dear \" run in the greeting line.\n//   2) Drops the <w:proofErr .../> markers around \"contact \" and merges\n//      it back into the same run as the preceding sentence (\"...hesitate\n//      to contact \").\n//   3) Collapses the four separate instrText runs that spelled out the\n//      ' SET SIGNATURE \"GARY\"' field code into a single run reading\n//      ' SET SIGNATURE GARY' (the stray quote-mark runs are removed).\n//\n// proofErr elements carry no text and aren't reachable through the\n// Word JS object model (Range/Paragraph text APIs never surface them),\n// so the simplest faithful way to drop them \u2014 and to let Word's own\n// OOXML serializer recombine runs that are now identical in\n// formatting/content \u2014 is a read-modify-write of each affected\n// paragraph's OOXML.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- Paragraph 1: \"Dear <<Greeting>>,\" -----------------------------------\n// Round-tripping through getOoxml()/insertOoxml() is enough: the shim's\n// OOXML reader never emits <w:proofErr/> markers in the first place, so\n// simply re-serializing the paragraph removes the gramStart/gramEnd pair\n// around the \"Dear \" run without touching anything else.\nconst greetingPara = paragraphs.items[0];\nconst greetingOoxml = greetingPara.getOoxml();\nawait context.sync();\ngreetingPara.insertOoxml(greetingOoxml.value, Word.InsertLocation.replace);\nawait context.sync();\n\n// --- \"...hesitate to contact <<ContactInfo>>.\" paragraph -----------------\n// Find it by its distinctive text rather than a hard-coded index, then\n// apply the same read-modify-write: the reader drops the proofErr pair\n// around \"contact \" and, because that run now shares identical rPr with\n// the run before it, Word's writer fuses them back into one <w:t> run \u2014\n// exactly matching the diff.\nconst contactSearch = context.document.body.search(\"hesitate to\", { matchCase: false });\ncontactSearch.load(\"items\");\nawait context.sync();\nconst contactPara = contactSearch.items[0].paragraphs.getFirst();\nconst contactOoxml = contactPara.getOoxml();\nawait context.sync();\ncontactPara.insertOoxml(contactOoxml.value, Word.InsertLocation.replace);\nawait context.sync();\n\n// --- \"Sincerely, SET SIGNATURE \"GARY\"\"  paragraph -------------------------\n// Locate the closing paragraph that holds the SET SIGNATURE field code\n// and collapse its four instrText runs (\" SET SIGNATURE \", '\"', \"GARY\",\n// '\"') into the single run ' SET SIGNATURE GARY' that the diff leaves\n// behind (the literal quote-mark runs are dropped, not just merged).\nconst signatureSearch = context.document.body.search(\"Sincerely\", { matchCase: false });\nsignatureSearch.load(\"items\");\nawait context.sync();\nconst signaturePara = signatureSearch.items[0].paragraphs.getFirst();\nconst signatureOoxmlResult = signaturePara.getOoxml();\nawait context.sync();\n\nlet signatureXml = signatureOoxmlResult.value;\nconst setSignatureRunsRe = /(<w:r>(<w:rPr>[\\s\\S]*?<\\/w:rPr>)<w:instrText xml:space=\"preserve\"> SET SIGNATURE <\\/w:instrText><\\/w:r>)<w:r><w:rPr>[\\s\\S]*?<\\/w:rPr><w:instrText>\"<\\/w:instrText><\\/w:r><w:r><w:rPr>[\\s\\S]*?<\\/w:rPr><w:instrText>GARY<\\/w:instrText><\\/w:r><w:r><w:rPr>[\\s\\S]*?<\\/w:rPr><w:instrText>\"<\\/w:instrText><\\/w:r>/;\nconst match = signatureXml.match(setSignatureRunsRe);\nif (match) {\n  const runProps = match[2];\n  const mergedRun = \"<w:r>\" + runProps + '<w:instrText xml:space=\"preserve\"> SET SIGNATURE GARY</w:instrText></w:r>';\n  signatureXml = signatureXml.replace(setSignatureRunsRe, mergedRun);\n  signaturePara.insertOoxml(signatureXml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Ver 7 - Mid Level Contact\n#\n# The commit being replayed:\n#   1) Drops the stray proofErr (gramStart/gramEnd) markers that bracket\n#      the \"Dear \" run in the greeting line.\n#   2) Drops the proofErr markers around \"contact \" and merges it back\n#      into the same run as the preceding sentence (\"...hesitate to\n#      contact \").\n#   3) Collapses the four separate instrText runs that spelled out the\n#      ' SET SIGNATURE \"GARY\"' field code into a single run reading\n#      ' SET SIGNATURE GARY' (the stray quote-mark runs are removed).\n#\n# proofErr elements carry no text and the Word object model's Find/Range\n# text never surfaces them, so the simplest faithful way to drop them --\n# and to let Word's own OOXML writer recombine runs that are now\n# identical in formatting/content -- is a read-modify-write of each\n# affected paragraph's Range.WordOpenXML via Range.InsertXML().\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphContaining($doc, [string]$needle) {\n    foreach ($para in $doc.Paragraphs) {\n        if ($para.Range.Text -match [regex]::Escape($needle)) {\n            return $para\n        }\n    }\n    return $null\n}\n\n# --- Paragraph 1: \"Dear <<Greeting>>,\" ------------------------------------\n# Round-tripping through WordOpenXML/InsertXML is enough: the reader never\n# emits <w:proofErr/> markers in the first place, so simply re-serializing\n# the paragraph removes the gramStart/gramEnd pair around the \"Dear \" run\n# without touching anything else.\n$greetingPara = Find-ParagraphContaining $d \"Dear \"\n$greetingRange = $greetingPara.Range\n$greetingRange.InsertXML($greetingRange.WordOpenXML) | Out-Null\n\n# --- \"...hesitate to contact <<ContactInfo>>.\" paragraph -------------------\n# Same read-modify-write: the reader drops the proofErr pair around\n# \"contact \" and, because that run now shares identical rPr with the run\n# before it, Word's writer fuses them back into one <w:t> run -- exactly\n# matching the diff.\n$contactPara = Find-ParagraphContaining $d \"hesitate to\"\n$contactRange = $contactPara.Range\n$contactRange.InsertXML($contactRange.WordOpenXML) | Out-Null\n\n# --- \"Sincerely, SET SIGNATURE \"GARY\"\"  paragraph ---------------------------\n# Locate the closing paragraph that holds the SET SIGNATURE field code and\n# collapse its four instrText runs (\" SET SIGNATURE \", '\"', \"GARY\", '\"')\n# into the single run ' SET SIGNATURE GARY' that the diff leaves behind\n# (the literal quote-mark runs are dropped, not just merged).\n$signaturePara = Find-ParagraphContaining $d \"Sincerely\"\n$signatureRange = $signaturePara.Range\n$signatureXml = $signatureRange.WordOpenXML\n\n$setSignatureRunsPattern = '(<w:r>(<w:rPr>[\\s\\S]*?</w:rPr>)<w:instrText xml:space=\"preserve\"> SET SIGNATURE </w:instrText></w:r>)<w:r><w:rPr>[\\s\\S]*?</w:rPr><w:instrText>\"</w:instrText></w:r><w:r><w:rPr>[\\s\\S]*?</w:rPr><w:instrText>GARY</w:instrText></w:r><w:r><w:rPr>[\\s\\S]*?</w:rPr><w:instrText>\"</w:instrText></w:r>'\n$match = [regex]::Match($signatureXml, $setSignatureRunsPattern)\nif ($match.Success) {\n    $runProps = $match.Groups[2].Value\n    $mergedRun = \"<w:r>\" + $runProps + '<w:instrText xml:space=\"preserve\"> SET SIGNATURE GARY</w:instrText></w:r>'\n    $signatureXml = $signatureXml -replace $setSignatureRunsPattern, [System.Text.RegularExpressions.Regex]::Replace($mergedRun, '\\$', '$$$$')\n    $signatureRange.InsertXML($signatureXml) | Out-Null\n}\n"}
